$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values that parse as plain numbers need a leading apostrophe so
# Excel stores them as text (matching the source data, which is always
# text in this sheet), not as a numeric value.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.080.90"
$ws.Range("E2").Value = "  +0.38%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.639.75"
$ws.Range("E3").Value = "  -0.07%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.34%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'214.54"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.27%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.42%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -2.43%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -2.24%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'18.58"
$ws.Range("E10").Value = "  -5.33%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0795"
$ws.Range("E11").Value = "  -0.21%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.744.82"
$ws.Range("E12").Value = "  +5.06%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  -1.63%  "

# Row 14 - Polygon
$ws.Range("E14").Value = "  -2.66%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "'62.37"
$ws.Range("E15").Value = "  -1.04%  "

# Row 16 & 17 - ShibaInu / WrappedBTC swap places
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "26.082.16"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0₃0748"
$ws.Range("E17").Value = "  -2.02%  "

# Row 18 - Dai
$ws.Range("E18").Value = "  +0.43%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'190.40"
$ws.Range("E19").Value = "  -1.41%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -2.09%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "'9.58"
$ws.Range("E21").Value = "  -3.53%  "

# Row 22 - Chainlink
$ws.Range("E22").Value = "  -2.51%  "

# Row 23 - Monero
$ws.Range("D23").Value = "'144.28"
$ws.Range("E23").Value = "  +0.52%  "

# Row 24 - Stellar
$ws.Range("E24").Value = "  -0.64%  "

# Row 25 - BinanceUSD
$ws.Range("E25").Value = "  +0.34%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "'1.76"
$ws.Range("E26").Value = "  -1.62%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'6.76"
$ws.Range("E27").Value = "  -1.53%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  -2.39%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -0.53%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -3.44%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  -2.40%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -3.61%  "

# Row 33 - HuobiToken
$ws.Range("D33").Value = "'2.44"
$ws.Range("E33").Value = "  -0.33%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "'1.51"
$ws.Range("E34").Value = "  -1.68%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  -2.47%  "

# Row 36 - Maker
$ws.Range("D36").Value = "1.121.74"
$ws.Range("E36").Value = "  -1.24%  "

# Row 37 - MXToken
$ws.Range("E37").Value = "  -0.20%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -3.94%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -1.64%  "

# Row 40 - Quant
$ws.Range("D40").Value = "'98.82"
$ws.Range("E40").Value = "  -0.59%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "'0.786"
$ws.Range("E41").Value = "  -1.49%  "

# Row 42 - FraxShare
$ws.Range("E42").Value = "  -3.57%  "

# Row 43 - BabyDogeCoin
$ws.Range("E43").Value = "  -0.48%  "

# Row 44 - Aave
$ws.Range("D44").Value = "'55.28"
$ws.Range("E44").Value = "  -2.52%  "

# Row 45 - Cronos
$ws.Range("E45").Value = "  -1.80%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  +0.09%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  -0.11%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "'7.65"
$ws.Range("E48").Value = "  -0.34%  "

# Row 49 - USDD
$ws.Range("E49").Value = "  +0.21%  "

# Row 50 - Algorand
$ws.Range("D50").Value = "'0.0930"
$ws.Range("E50").Value = "  -3.36%  "

# Row 51 - NEARProtocol
$ws.Range("E51").Value = "  -1.18%  "
